$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Remove the trailing two "Policy Subsystems" slides (was slide38.xml /
#    slide39.xml, the last two slides in the deck).
# ---------------------------------------------------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
$p.Slides.Item($lastIndex - 1).Delete()

# ---------------------------------------------------------------------------
# 2) Title slide: the static "date" placeholder text changes from the
#    generic "1/1/23" to "Spring 2023".
# ---------------------------------------------------------------------------
$titleSlide = $p.Slides.Item(1)
for ($i = 1; $i -le $titleSlide.Shapes.Count; $i++) {
    $shp = $titleSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "1/1/23") {
        $shp.TextFrame.TextRange.Text = "Spring 2023"
    }
}

# ---------------------------------------------------------------------------
# 3) Every inherited "datetimeFigureOut" footer field (slide master + every
#    custom layout + notes master) gets re-cached from 2/13/23 -> 2/19/23.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $newText) {
    $placeholders = $shapes.Placeholders
    for ($j = 1; $j -le $placeholders.Count; $j++) {
        $shp = $placeholders.Item($j)
        if ($shp.PlaceholderFormat.Type -eq 16 -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -ne "") {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes "2/19/23"

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes "2/19/23"
}

Update-DatePlaceholder $p.NotesMaster.Shapes "2/19/23"
